$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of seasonal stats for "M2_09 Dryad 2020",
# copying the formatting (border/alignment) of the prior season's row.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "M2_09 Dryad 2020"
$ws.Range("C10").Value = 9678
$ws.Range("D10").Value = 10725
$ws.Range("E10").Value = 855528
$ws.Range("F10").Value = 9946
$ws.Range("G10").Value = 10046
$ws.Range("H10").Value = 10183
